$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the old "Overall Mass (ABS)" totals row (row 48).
#    This also drops the stray dependency on H2 (density value).
# ------------------------------------------------------------------
$ws.Rows.Item(48).Delete()

# ------------------------------------------------------------------
# 2. Add the new "Material/Process" column header in G1, matching
#    the bold / filled header formatting already used by B1/D1.
# ------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Material/Process"

# ------------------------------------------------------------------
# 3. Replace the old density mini-table in G2:I2
#    ("(ABS) Density" / 1.06 / "g/cm^3") with a per-part
#    Material/Process tag in column G, and clear H2/I2.
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("G2:G45").PasteSpecial(-4122)

$ws.Range("G2:G45").Value = "ASA/FDM"
$ws.Range("G5:G8").Value  = "Nylon/SLS"
$ws.Range("G18").Value    = "Nylon/SLS"
$ws.Range("G26").Value    = "Nylon/SLS"
$ws.Range("G42").Value    = "Nylon/SLS"
$ws.Range("G45").Value    = "Nylon/SLS"

$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()

Write-Host "done"
